# Update the "Last Updated" timestamp on the Metadata sheet
$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "29 Oct 2025, 04:01 PM"

# Update the "Top Losers" sheet: a new MPSLTD entry was inserted at row 33,
# shifting the previous rows 33-38 down to 34-39 (row 39's old MPSLTD entry
# is dropped, row 40 BOSCHLTD stays as-is).
$wsLosers = $wb.Worksheets.Item("Top Losers")

# New row inserted at 33
$wsLosers.Range("B33").Value = "MPSLTD"
$wsLosers.Range("C33").Value = -3.2782
$wsLosers.Range("D33").Value = -4.6314
$wsLosers.Range("E33").Value = 2.1755

# Row 34 <- old row 33 (CAMS)
$wsLosers.Range("B34").Value = "CAMS"
$wsLosers.Range("C34").Value = -3.2545
$wsLosers.Range("D34").Value = -0.6366000000000001
$wsLosers.Range("E34").Value = 2.5781

# Row 35 <- old row 34 (PRUDENT)
$wsLosers.Range("B35").Value = "PRUDENT"
$wsLosers.Range("C35").Value = -3.2484
$wsLosers.Range("D35").Value = -3.6312
$wsLosers.Range("E35").Value = 1.9933

# Row 36 <- old row 35 (SPARC)
$wsLosers.Range("B36").Value = "SPARC"
$wsLosers.Range("C36").Value = -3.1709
$wsLosers.Range("D36").Value = 4.8337
$wsLosers.Range("E36").Value = 6.3311

# Row 37 <- old row 36 (ANANDRATHI)
$wsLosers.Range("B37").Value = "ANANDRATHI"
$wsLosers.Range("C37").Value = -3.0775
$wsLosers.Range("D37").Value = -0.8672
$wsLosers.Range("E37").Value = 9.1835

# Row 38 <- old row 37 (NLCINDIA)
$wsLosers.Range("B38").Value = "NLCINDIA"
$wsLosers.Range("C38").Value = -3.0757
$wsLosers.Range("D38").Value = -4.5618
$wsLosers.Range("E38").Value = -11.6431

# Row 39 <- old row 38 (YATRA)
$wsLosers.Range("B39").Value = "YATRA"
$wsLosers.Range("C39").Value = -3.0403
$wsLosers.Range("D39").Value = -2.8455
$wsLosers.Range("E39").Value = 7.3711
